$wb = $excel.ActiveWorkbook

# --- Sheet "Info": update B3 value ---
$infoWs = $wb.Worksheets.Item("Info")
$infoWs.Range("B3").Value = 0.754

# --- Sheet "Solution": replace numeric codes with their text labels ---
# Values are written column-by-column (B, C, D, E, F) top-to-bottom so that
# the resulting shared-string table is built up in the same order as the
# target workbook.
$solWs = $wb.Worksheets.Item("Solution")

$solWs.Range("B17").Value = "FEP 0.0625"

$solWs.Range("C9").Value  = "CXM 2.25"
$solWs.Range("C15").Value = "AMC 8_8"
$solWs.Range("C16").Value = "CXM 2.25"

$solWs.Range("D5").Value  = "ZOX 0.0156"
$solWs.Range("D7").Value  = "CEC 4"
$solWs.Range("D8").Value  = "CPR 16"
$solWs.Range("D9").Value  = "CXM 2.25"
$solWs.Range("D11").Value = "CEC 4"
$solWs.Range("D12").Value = "AM 512"
$solWs.Range("D14").Value = "TZP 8_32"
$solWs.Range("D16").Value = "CXM 2.25"
$solWs.Range("D17").Value = "I"

$solWs.Range("E3").Value  = "AMC 8_8"
$solWs.Range("E4").Value  = "CPR 12"
$solWs.Range("E5").Value  = "ZOX 0.0156"
$solWs.Range("E6").Value  = "I"
$solWs.Range("E7").Value  = "CEC 4"
$solWs.Range("E8").Value  = "FEP 0.0312"
$solWs.Range("E9").Value  = "I"
$solWs.Range("E10").Value = "AMC 8_8"
$solWs.Range("E11").Value = "CEC 4"
$solWs.Range("E12").Value = "AM 512"
$solWs.Range("E13").Value = "I"
$solWs.Range("E14").Value = "TZP 8_32"
$solWs.Range("E15").Value = "I"
$solWs.Range("E16").Value = "I"
$solWs.Range("E17").Value = "I"

$solWs.Range("F2").Value  = "I"
$solWs.Range("F3").Value  = "AMC 8_8"
$solWs.Range("F4").Value  = "AM 512"
$solWs.Range("F5").Value  = "I"
$solWs.Range("F6").Value  = "AMC 8_8"
$solWs.Range("F9").Value  = "I"
$solWs.Range("F10").Value = "CPR 12"
$solWs.Range("F11").Value = "I"
$solWs.Range("F13").Value = "I"
$solWs.Range("F15").Value = "I"
$solWs.Range("F16").Value = "I"
$solWs.Range("F17").Value = "I"
